$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 'Running Robot tests (Assigning) studying Assigning generously Assign all operators Assign all operators for the remaining users'
$ws.Range("F2").Value = '[''Running'', ''Robot'', ''tests'', ''('', ''Assigning'', '')'', ''studying'', ''Assigning'', ''generously'', ''Assign'', ''all'', ''operators'', ''Assign'', ''all'', ''operators'', ''for'', ''the'', ''remaining'', ''users'']'
$ws.Range("G2").Value = '[''Running'', ''Robot'', ''tests'', ''Assigning'', ''studying'', ''Assigning'', ''generously'', ''Assign'', ''all'', ''operators'', ''Assign'', ''all'', ''operators'', ''for'', ''the'', ''remaining'', ''users'']'
$ws.Range("E15").Value = 'Reports tests (Cleaning Jobs) I Me Verify "Download Report" button is disabled by default and the empty PDF message in the container is shown Input worksite, robot and month, then select first report shown, assert API call, download the report then verify the report has been downloaded (Excluding CC) I Me Input company, worksite, robot and month but do not select any report, verify "Download Report" button is disabled and correct empty PDF frame (Excluding CC) Input company, worksite, robot and month, download the first report, verify successful api call and file download user 32131 use213123r 5434634'
$ws.Range("F15").Value = '[''Reports'', ''tests'', ''('', ''Cleaning'', ''Jobs'', '')'', ''I'', ''Me'', ''Verify'', ''``'', ''Download'', ''Report'', "''''", ''button'', ''is'', ''disabled'', ''by'', ''default'', ''and'', ''the'', ''empty'', ''PDF'', ''message'', ''in'', ''the'', ''container'', ''is'', ''shown'', ''Input'', ''worksite'', '','', ''robot'', ''and'', ''month'', '','', ''then'', ''select'', ''first'', ''report'', ''shown'', '','', ''assert'', ''API'', ''call'', '','', ''download'', ''the'', ''report'', ''then'', ''verify'', ''the'', ''report'', ''has'', ''been'', ''downloaded'', ''('', ''Excluding'', ''CC'', '')'', ''I'', ''Me'', ''Input'', ''company'', '','', ''worksite'', '','', ''robot'', ''and'', ''month'', ''but'', ''do'', ''not'', ''select'', ''any'', ''report'', '','', ''verify'', ''``'', ''Download'', ''Report'', "''''", ''button'', ''is'', ''disabled'', ''and'', ''correct'', ''empty'', ''PDF'', ''frame'', ''('', ''Excluding'', ''CC'', '')'', ''Input'', ''company'', '','', ''worksite'', '','', ''robot'', ''and'', ''month'', '','', ''download'', ''the'', ''first'', ''report'', '','', ''verify'', ''successful'', ''api'', ''call'', ''and'', ''file'', ''download'', ''user'', ''32131'', ''use213123r'', ''5434634'']'
$ws.Range("G15").Value = '[''Reports'', ''tests'', ''Cleaning'', ''Jobs'', ''I'', ''Me'', ''Verify'', ''Download'', ''Report'', ''button'', ''is'', ''disabled'', ''by'', ''default'', ''and'', ''the'', ''empty'', ''PDF'', ''message'', ''in'', ''the'', ''container'', ''is'', ''shown'', ''Input'', ''worksite'', ''robot'', ''and'', ''month'', ''then'', ''select'', ''first'', ''report'', ''shown'', ''assert'', ''API'', ''call'', ''download'', ''the'', ''report'', ''then'', ''verify'', ''the'', ''report'', ''has'', ''been'', ''downloaded'', ''Excluding'', ''CC'', ''I'', ''Me'', ''Input'', ''company'', ''worksite'', ''robot'', ''and'', ''month'', ''but'', ''do'', ''not'', ''select'', ''any'', ''report'', ''verify'', ''Download'', ''Report'', ''button'', ''is'', ''disabled'', ''and'', ''correct'', ''empty'', ''PDF'', ''frame'', ''Excluding'', ''CC'', ''Input'', ''company'', ''worksite'', ''robot'', ''and'', ''month'', ''download'', ''the'', ''first'', ''report'', ''verify'', ''successful'', ''api'', ''call'', ''and'', ''file'', ''download'', ''user'', ''32131'', ''use213123r'', ''5434634'']'
